$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 147, pushing existing rows 147:230 down to 148:231
$ws.Rows.Item(147).Insert()

# Populate the newly inserted row 147 with the new record's data
$ws.Cells.Item(147, 1).Value = 10
$ws.Cells.Item(147, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(147, 3).Value = "La Araucanía"
$ws.Cells.Item(147, 4).Value = 44455
$ws.Cells.Item(147, 5).Value = 9
$ws.Cells.Item(147, 6).Value = 100112032
$ws.Cells.Item(147, 7).Value = "Zapallo italiano"
$ws.Cells.Item(147, 8).Value = "Sin especificar"
$ws.Cells.Item(147, 9).Value = "Primera"
$ws.Cells.Item(147, 10).Value = 70
$ws.Cells.Item(147, 11).Value = 18000
$ws.Cells.Item(147, 12).Value = 19000
$ws.Cells.Item(147, 13).Value = 18429
$ws.Cells.Item(147, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(147, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(147, 16).Value = 307
$ws.Cells.Item(147, 17).Value = 60
$ws.Cells.Item(147, 18).Value = "Hortaliza"
